$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, pushing existing rows 233:340 down to 234:341
$ws.Rows(233).Insert()

# Populate the new row 233 with the weekly data point that was added
$ws.Range("A233").Value = 7
$ws.Range("B233").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C233").Value = "Ñuble"
$ws.Range("D233").Value = (Get-Date -Year 2023 -Month 6 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E233").Value = 16
$ws.Range("F233").Value = 100112043
$ws.Range("G233").Value = "Pepino ensalada"
$ws.Range("H233").Value = "Sin especificar"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 120
$ws.Range("K233").Value = 12000
$ws.Range("L233").Value = 13000
$ws.Range("M233").Value = 12500
$ws.Range("N233").Value = "$/caja 60 unidades"
$ws.Range("O233").Value = "Región de Arica y Parinacota"
$ws.Range("P233").Value = 208
$ws.Range("Q233").Value = 60
$ws.Range("R233").Value = "Hortaliza"
